$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.96
$ws.Range("K2").Value = 4
$ws.Range("P2").Value = 1.92
$ws.Range("W2").Value = 2.04
$ws.Range("G3").Value = 2.68
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 4.1
$ws.Range("K3").Value = 3.5
$ws.Range("L3").Value = 1.48
$ws.Range("N3").Value = 2.92
$ws.Range("O3").Value = 1.42
$ws.Range("Q3").Value = 2.24
$ws.Range("S3").Value = 4.3
$ws.Range("U3").Value = 1.9
$ws.Range("V3").Value = 1.33
$ws.Range("X3").Value = 13
$ws.Range("Y3").Value = 13.5
$ws.Range("AB3").Value = 10.5
$ws.Range("O4").Value = 1.23
$ws.Range("P4").Value = 2.18
$ws.Range("Q4").Value = 1.67
$ws.Range("T4").Value = 2.84
$ws.Range("U4").Value = 1.44
$ws.Range("AB4").Value = 9
$ws.Range("F5").Value = 1.76
$ws.Range("J5").Value = 3.85
$ws.Range("M5").Value = 1.03
$ws.Range("P5").Value = 2.26
$ws.Range("S5").Value = 2.56
$ws.Range("T5").Value = 1.63
$ws.Range("AL5").Value = 30
$ws.Range("L7").Value = 1.25
$ws.Range("N7").Value = 5
$ws.Range("S7").Value = 2.48
$ws.Range("AN7").Value = 12
$ws.Range("F8").Value = 2.22
$ws.Range("G8").Value = 2.52
$ws.Range("H8").Value = 2.86
$ws.Range("I8").Value = 3.15
$ws.Range("S8").Value = 1.94
$ws.Range("T8").Value = 1.4
$ws.Range("V8").Value = 1.46
$ws.Range("W8").Value = 1.66
$ws.Range("AN8").Value = 11
$ws.Range("S9").Value = 1.99
$ws.Range("F10").Value = 1.8
$ws.Range("L10").Value = 1.2
$ws.Range("Q10").Value = 1.42
$ws.Range("S10").Value = 2.06
$ws.Range("H11").Value = 3.8
$ws.Range("J11").Value = 3.2
$ws.Range("N11").Value = 3.15
$ws.Range("Q11").Value = 2.08
$ws.Range("AC11").Value = 9.4
$ws.Range("AK11").Value = 30
$ws.Range("AM11").Value = 150
$ws.Range("AN11").Value = 23
$ws.Range("AO11").Value = 85
$ws.Range("H12").Value = 2.74
$ws.Range("M12").Value = 1.11
$ws.Range("U12").Value = 1.9
$ws.Range("W12").Value = 1.44
$ws.Range("Z12").Value = 18
$ws.Range("AA12").Value = 55
$ws.Range("AD12").Value = 13.5
$ws.Range("AK12").Value = 50
$ws.Range("I13").Value = 16
